# Adds a new record (Paula / paula_@hotmail.com / 123) as row 10 on Hoja1,
# turning the email cell into a live mailto: hyperlink, and moves the
# active selection to C11 (just past the new last row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New data row
$ws.Range("A10").Value = "Paula"
$ws.Range("B10").Value = "paula_@hotmail.com"
$ws.Range("C10").Value = 123

# Turn the email into a hyperlink (mailto:), which also applies Excel's
# built-in "Hyperlink" cell style (underline + theme color font) to B10.
$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:paula_@hotmail.com") | Out-Null

# Match the saved selection state from the source workbook.
$ws.Range("C11").Select() | Out-Null
